$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Helper donor cells (row 23 is untouched by this edit, so it is safe to use
# throughout as a style donor via copy / paste-special-formats):
#   C23 / D23 / G23  -> style "14" text cell holding "0"      (shared idx 20)
#   E23 / H23        -> style "14" text cell holding "***.*"  (shared idx 21)
#   F23               -> style "16" numeric cell
#   K23               -> style "15" numeric cell
# ---------------------------------------------------------------------------

function Set-TextCell($addr, $donor, $text) {
    # Force the destination to text first (apostrophe prefix => text value),
    # then paste just the donor's formatting on top, which both applies the
    # desired style AND clears the quote-prefix flag that the text entry set.
    $ws.Range($addr).Value = "'" + $text
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
}

function Set-NumberCell($addr, $donor, $num) {
    $ws.Range($donor).Copy()
    $ws.Range($addr).PasteSpecial(-4122)
    $ws.Range($addr).Value = $num
}

# --- sharedStrings text edits -----------------------------------------------
$ws.Range("A8").Value = "Volume 31   Number  30"
$ws.Range("C9").Value = "Report Covering the Week  7/22/2024  Through  7/28/2024"

# --- Row 15 (Rape) ------------------------------------------------------
Set-TextCell "G15" "C23" "0"
Set-TextCell "H15" "E23" "***.*"

# --- Row 16 (Robbery) ----------------------------------------------------
Set-NumberCell "D16" "F23" 2
Set-NumberCell "E16" "K23" 0
$ws.Range("G16").Value = 5
$ws.Range("H16").Value = 40
$ws.Range("I16").Value = 48
$ws.Range("J16").Value = 39
$ws.Range("K16").Value = 23.076923076923
$ws.Range("L16").Value = 17.073170731707
$ws.Range("M16").Value = -15.78947368421
$ws.Range("N16").Value = -87.909319899244

# --- Row 17 (Fel. Assault) -------------------------------------------------
$ws.Range("C17").Value = 4
Set-TextCell "D17" "C23" "0"
Set-TextCell "E17" "E23" "***.*"
$ws.Range("F17").Value = 12
$ws.Range("G17").Value = 6
$ws.Range("H17").Value = 100
$ws.Range("I17").Value = 62
$ws.Range("K17").Value = 8.771929824561
$ws.Range("L17").Value = 10.714285714285
$ws.Range("M17").Value = 82.35294117647
$ws.Range("N17").Value = -10.144927536231

# --- Row 18 (Burglary) ------------------------------------------------------
Set-TextCell "C18" "C23" "0"
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = -100
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -50
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = -30
$ws.Range("L18").Value = -20.967741935483
$ws.Range("M18").Value = -34.666666666666
$ws.Range("N18").Value = -93.518518518518

# --- Row 19 (Gr. Larceny) ---------------------------------------------------
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -55.555555555555
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 42
$ws.Range("H19").Value = -26.190476190476
$ws.Range("I19").Value = 254
$ws.Range("J19").Value = 284
$ws.Range("K19").Value = -10.56338028169
$ws.Range("L19").Value = -12.714776632302
$ws.Range("M19").Value = 22.705314009661
$ws.Range("N19").Value = -54.151624548736

# --- Row 20 (G.L.A.) --------------------------------------------------------
$ws.Range("C20").Value = 4
Set-TextCell "D20" "C23" "0"
Set-TextCell "E20" "E23" "***.*"
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 21
$ws.Range("H20").Value = -57.142857142857
$ws.Range("I20").Value = 96
$ws.Range("K20").Value = 14.285714285714
$ws.Range("L20").Value = 84.615384615384
$ws.Range("M20").Value = 74.545454545454
$ws.Range("N20").Value = -95.131845841785

# --- Row 21 (TOTAL) ----------------------------------------------------------
$ws.Range("C21").Value = 14
$ws.Range("D21").Value = 13
$ws.Range("E21").Value = 7.692307692307
$ws.Range("F21").Value = 63
$ws.Range("G21").Value = 82
$ws.Range("H21").Value = -23.170731707317
$ws.Range("I21").Value = 514
$ws.Range("J21").Value = 537
$ws.Range("K21").Value = -4.283054003724
$ws.Range("L21").Value = 0.587084148727
$ws.Range("M21").Value = 19.53488372093
$ws.Range("N21").Value = -86.311584553928

# --- Row 22 (Transit) --------------------------------------------------------
Set-TextCell "C22" "C23" "0"
Set-TextCell "D22" "C23" "0"
Set-TextCell "E22" "E23" "***.*"
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = -50

# --- Row 24 (Petit Larceny) --------------------------------------------------
$ws.Range("C24").Value = 39
$ws.Range("D24").Value = 19
$ws.Range("E24").Value = 105.263157894737
$ws.Range("F24").Value = 126
$ws.Range("G24").Value = 89
$ws.Range("H24").Value = 41.573033707865
$ws.Range("I24").Value = 946
$ws.Range("J24").Value = 874
$ws.Range("K24").Value = 8.237986270022
$ws.Range("L24").Value = -6.336633663366
$ws.Range("M24").Value = 68.627450980392

# --- Row 25 (Retail Theft) ---------------------------------------------------
$ws.Range("C25").Value = 34
$ws.Range("D25").Value = 12
$ws.Range("E25").Value = 183.333333333333
$ws.Range("F25").Value = 95
$ws.Range("G25").Value = 58
$ws.Range("H25").Value = 63.793103448275
$ws.Range("I25").Value = 687
$ws.Range("J25").Value = 622
$ws.Range("K25").Value = 10.450160771704
$ws.Range("L25").Value = -4.18410041841

# --- Row 26 (Misd. Assault) ---------------------------------------------------
$ws.Range("C26").Value = 4
$ws.Range("D26").Value = 3
$ws.Range("E26").Value = 33.333333333333
$ws.Range("F26").Value = 21
$ws.Range("G26").Value = 17
$ws.Range("H26").Value = 23.529411764705
$ws.Range("I26").Value = 163
$ws.Range("J26").Value = 136
$ws.Range("K26").Value = 19.85294117647
$ws.Range("L26").Value = 53.77358490566
$ws.Range("M26").Value = 27.34375

# --- Row 27 (UCR Rape*) ---------------------------------------------------
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 0

# --- Row 28 (Other Sex Crimes) -----------------------------------------------
$ws.Range("D28").Value = 4
Set-TextCell "F28" "C23" "0"
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = -100
$ws.Range("J28").Value = 20
$ws.Range("K28").Value = -10
$ws.Range("L28").Value = -14.285714285714
